# Updates the cryptos price/volume tracking sheet (auto-refreshed data).
# Column D = Price (text), Column E = Volume(1h) (text, padded with spaces).
# Values are entered with a leading apostrophe to force text (so Excel does
# not reinterpret numeric-looking prices like "142.00" as a Number), then the
# quote-prefix cell style introduced by that is cleared by resetting to Normal
# so the saved cell style matches the original (unstyled) data cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.091.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.32%  "

# Row 3
$ws.Range("D3").Value = "'2.894.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'525.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.38%  "

# Row 6
$ws.Range("D6").Value = "'142.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.69%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.549"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.47%  "

# Row 9
$ws.Range("D9").Value = "'2.891.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.91%  "

# Row 10
$ws.Range("D10").Value = "'0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.94%  "

# Row 11
$ws.Range("E11").Value = "  -7.89%  "

# Row 12
$ws.Range("D12").Value = "'0.354"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.70%  "

# Row 13
$ws.Range("D13").Value = "'3.398.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.64%  "

# Row 14
$ws.Range("E14").Value = "  +1.16%  "

# Row 15
$ws.Range("D15").Value = "'60.279.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.14%  "

# Row 16
$ws.Range("D16").Value = "'22.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.32%  "

# Row 17
$ws.Range("D17").Value = "'2.893.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.88%  "

# Row 18
$ws.Range("E18").Value = "  -6.30%  "

# Row 19
$ws.Range("D19").Value = "'4.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.13%  "

# Row 20
$ws.Range("D20").Value = "'11.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.94%  "

# Row 21
$ws.Range("D21").Value = "'359.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.93%  "

# Row 22
$ws.Range("D22").Value = "'6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "'5.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").Value = "'63.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.90%  "

# Row 26
$ws.Range("D26").Value = "'3.017.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "

# Row 27
$ws.Range("D27").Value = "'0.448"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "

# Row 28
$ws.Range("E28").Value = "  -7.89%  "

# Row 29
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("D30").Value = "'7.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.92%  "

# Row 31
$ws.Range("D31").Value = "'0.0₃0847"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.77%  "

# Row 32
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("D33").Value = "'1.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.64%  "

# Row 34
$ws.Range("D34").Value = "'19.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.13%  "

# Row 35
$ws.Range("D35").Value = "'151.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.94%  "

# Row 36
$ws.Range("D36").Value = "'4.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.47%  "

# Row 37
$ws.Range("D37").Value = "'5.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.04%  "

# Row 38
$ws.Range("D38").Value = "'0.978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.80%  "

# Row 39
$ws.Range("D39").Value = "'1.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.61%  "

# Row 40
$ws.Range("D40").Value = "'37.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("D41").Value = "'2.332.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.00%  "

# Row 42
$ws.Range("D42").Value = "'1.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.76%  "

# Row 43 - was Mantle, now Filecoin (rows 43 and 44 swapped order)
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.26%  "

# Row 44 - was Filecoin, now Mantle
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.641"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.29%  "

# Row 45
$ws.Range("D45").Value = "'20.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.85%  "

# Row 46
$ws.Range("D46").Value = "'0.0566"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.51%  "

# Row 47
$ws.Range("D47").Value = "'0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "

# Row 48
$ws.Range("D48").Value = "'4.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.34%  "

# Row 49
$ws.Range("E49").Value = "  -1.64%  "

# Row 50
$ws.Range("D50").Value = "'0.0232"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.71%  "

# Row 51
$ws.Range("D51").Value = "'0.0925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
